# "Site updated" rebuild: refresh the total-time counters on the times
# sheet (hours / mins / secs of total site reading/writing time).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The counters are stored as text (they trip Excel's "number stored as
# text" warning in the original file), so force the Text format before
# writing the new values to keep them as text instead of being
# auto-coerced to numbers.
$counters = $ws.Range("A2:C2")
$counters.NumberFormat = "@"

$ws.Range("A2").Value = "118"
$ws.Range("B2").Value = "44"
$ws.Range("C2").Value = "37"
